$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.98%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.85"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "1.94%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.084"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.76%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08029"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.11%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.912"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.36%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.51%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.732"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.86%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9274"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.71%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1368"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.05%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.43%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09158"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.66%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03595"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.19%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09821"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.25%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001417"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.78%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005913"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.83%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.558"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.51%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.978"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2.22%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3455"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.86%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1305"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.53%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.909"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-2.62%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2511"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.76%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04461"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.01%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001225"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.73%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004791"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.12%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001486"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "18.71%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003136"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "4.24%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01957"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.65%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04909"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "4.41%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007649"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.72%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009153"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-6.19%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1374"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.94%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002108"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.21%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01136"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "18.22%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006399"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.72%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.09%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "63.57"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.41%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-19.89%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.09%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.09%"
